# Week 13 logging update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 315
$ws.Range("F2").Value = 114
$ws.Range("G2").Value = 78
$ws.Range("H2").Value = 12
$ws.Range("I2").Value = 12
$ws.Range("N2").Value = 39
$ws.Range("O2").Value = 40
$ws.Range("P2").Value = 23
$ws.Range("C3").Value = 314
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 58
$ws.Range("F3").Value = 182
$ws.Range("H3").Value = 60
$ws.Range("I3").Value = 113
$ws.Range("J3").Value = 92
$ws.Range("L3").Value = 497
$ws.Range("M3").Value = 309
$ws.Range("Q3").Value = 858

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 353
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = 116
$ws.Range("G2").Value = 103
$ws.Range("J2").Value = 54
$ws.Range("N2").Value = 28
$ws.Range("C3").Value = 318
$ws.Range("E3").Value = 54
$ws.Range("F3").Value = 203
$ws.Range("G3").Value = 58
$ws.Range("H3").Value = 29
$ws.Range("I3").Value = 112
$ws.Range("J3").Value = 100
$ws.Range("L3").Value = 513
$ws.Range("M3").Value = 349
$ws.Range("Q3").Value = 909

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 125
$ws.Range("D2").Value = 113
$ws.Range("J2").Value = 175
$ws.Range("K2").Value = 158
$ws.Range("N2").Value = 50
$ws.Range("O2").Value = 31
$ws.Range("B3").Value = 55

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B3").Value = 12
$ws.Range("E3").Value = 14

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("D2").Value = 11

$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = '-1 7 -5 1 -4 1 -2 -2 -2 5 -4 4 1 4 6 2 7 -1 13 -1 5 18 0 6 3 6 3 2 7 1 -1 3 2 12 6 0 4 0 6 19 3 2 1 23 -3 7 3 16 1 0 3 8 -1 3 4 4 4 2 4 4 10 3 3 1 3 6 5 10 26 4 3 3 0 13 11 9 2 9 3 2 3 8 2 4 0 1 3 2 5 1 2 4 1 7 -6 6 13 1 6 -2 4 14 -3 -3 5 4 49 0 -1 9 5 8 9 -1 1 8 3 4 0 3 7 2 7 1 4 7 7 10 3 2 3 2 -1 80 1 -2 2 4 0 19 -4 1 14 10 3 0 2 -2 2 2 1 9 1 6 11 3 17 6 4 2 2 1 2 6 5 4 -1 15 7 18 5 6 4 2 4 2 5 2 19 11 5 5 8 1 0 12 3 5 5 3 20 1 -2 2 2 6 -2 4 6 17 2 34 2 6 2 4 10 5 4 0 0 1 6 3 3 4 0 8 1 1 10 1 9 6 1 0 4 -3 6 2 -4 5 4 2 1 0 4 5 3 13 2 1 4 -1 2 4 5 18 4 7 2 4 2 5 2 2 -1 4 7 5 -2 3 -1 3 3 4 2 8 2 5 4 1 -1 6 4 5 3 3 0 5 9 2 3 60 13 4 3 13 23 3 4 0 3 1 5 3 6 -1 8 4 1 3 5 3 0 6 0 13 3 16 3 1 9 -1 4 9 2 2 5 -3 7 8 2 0 2 12 4 4 4 3 1 -3 1 14 3 5 6 1 1 9 3 3 3 15 3 8 10 7 23 3 0 1 1 2 11 4 0 3 5 2 2 7 0 5 11 -1 1 4 5 8 5 -2 4 2 -1 1 5 2 3 7 3 5 2 4 2 6 1 4 5 2 15 -2 2 6 41 6 2 7 5 3 0 9 -1 -6 2 -9 7 2 2 5 11 1 1 2 8 7 2 5 12 -3 15 2 -1 3 4 0 4 6 -3 2 4 2 6 11 4 4 1 1 6 -1 0 4 0 1 7 8 1 17 20 1 1 -3 2 4 1 4 2 1 6 6 4 6 5 1 4 3 4 1 14 2 4 1 0 0 -1 4 0 2 4 -2 7 3 4 4 14 9 2 3 0 2 2 0 6 -1 2 4 2 2 5 4 14 -1 7 2 2 0 1 9 3 0 -2 -2 7 11 3 4 2 3 6 2 2 -2 4 19 0 4 9 2 0 2 -3 11 2 2 6 2 14 5 4 0 3 5 4 10 2 2 6 4 16 5 5 6 2 -2 2 6 1 6 20 5 2 3 5 2 0 5 7 5 2 3 8 -3 4 11 12 2 2 8 6 0 1 3 16 1 2 12 0 9 6 5 0 2 14 2 2 1 1 32 6 -1 -3 7 -4 3 4 13 -1 1 4 3 5 0 1 -1 3 3 10 -2 0 1 0 7 0 23 17 1 0 9 13 3'
$wsYDS.Range("B3").Value = '13 14 7 1 11 12 41 14 18 5 13 10 4 4 5 5 4 15 10 6 7 7 5 7 5 17 15 24 12 13 -3 7 7 11 20 12 6 4 2 12 10 11 8 22 4 9 4 2 18 9 13 23 0 11 2 9 8 5 5 10 29 4 22 -7 8 17 -2 10 12 4 6 11 -1 8 12 8 9 7 -1 9 15 13 4 33 5 12 17 18 7 9 3 11 21 16 12 16 -3 7 10 39 3 9 14 -1 2 21 8 10 23 5 9 -4 1 9 18 1 11 39 5 7 -4 17 29 9 12 0 5 0 7 11 7 7 12 5 9 2 8 8 18 0 7 11 5 12 14 11 11 12 6 15 4 11 9 1 12 1 6 10 19 6 6 3 50 4 2 7 8 21 8 10 16 8 7 3 7 12 4 8 7 9 6 0 22 4 11 9 12 3 8 4 17 11 1 4 6 27 38 7 16 7 0 40 -3 6 4 53 6 12 9 11 7 10 18 16 5 7 5 44 2 -9 13 -1 8 15 5 9 22 16 14 11 3 7 -2 3 6 8 6 14 13 6 7 11 3 4 10 3 39 7 24 4 22 35 -1 12 9 14 0 7 1 21 8 10 14 13 24 5 9 12 6 10 7 14 13 9 11 11 4 15 16 19 20 8 -5 7 0 9 16 13 11 3 12 14 15 14 16 9 15 21 10 10 4 2 2 18 33 7 16 3 9 21 13 10 42 -6 4 17 1 6 8 37 2 15 20 8 17 14 6 16 16 6 9 15 8 6 12 16 13 4 4 8 2 6 17 3 4 2 8 12 19 33 13 21 10 8 21 13 14 8 38 19 8 18 12 9 7 4 16 0 13 13 11 10 12 3 7 5 16 4 9 10 13 -4 6 18 4 27 51 19 7 8 17 20 8 4 5 21 1 54 8 5 11 3 28 17 8 18 1 23 38 5 9 18 28 6 35 6 2 26 5 13 13 16 14 13 13 11 7 13 17 8 3 16 2 9 4 14 7 8 3 10 10 11 6 10 3 10 1 11 12 11 11 7 6 2 15 -2 6 18 7 1 27 15 7 16 15 1 6 8 15 -4 2 14 16 2 5 8 9 12 5 50 19 3 11 1 12 41 2 5 4 22 4 6 11 6 5 7 16 4 3 3 9 5 30 6 12 7 -7 14 3 14 2 9 1 -1 6 9 12 9 4 3 5 5 0 11 21 4 28 8 11 8 -1 7 2 4 5 5 1 20 7 18 19 6 14 9 0 10 17 17 1 18 18 5 4 10 8 3 15 -3 7 -5 4 18 20 -1 5 12 16 6 11 17 8 3 9 8 2 10 11 11'
$wsYDS.Range("C2").Value = '2 3 0 3 1 3 6 -1 1 8 1 5 11 3 14 5 0 21 -2 1 5 0 30 2 4 1 7 1 3 2 5 -1 4 7 1 1 3 8 10 1 3 3 1 2 0 4 4 3 2 2 5 6 11 10 -1 12 23 0 1 0 3 1 9 0 3 4 4 -3 -1 5 4 4 0 -1 -3 11 10 5 2 0 19 19 1 -1 2 3 1 23 3 6 3 -2 2 2 2 4 3 3 2 4 2 9 -1 13 3 4 0 3 5 1 -2 2 3 2 -2 2 9 9 5 5 3 8 3 8 1 0 1 9 7 1 13 9 12 0 0 7 4 7 -1 3 -2 2 4 1 14 5 0 0 0 6 0 0 3 8 -1 2 4 11 6 2 6 1 4 2 3 1 12 5 8 1 0 1 1 2 3 7 5 3 10 3 2 0 2 1 2 0 4 6 12 4 2 4 5 5 2 13 1 2 2 7 2 3 7 7 0 3 3 2 -2 2 7 0 2 2 5 1 10 -1 6 11 1 5 9 11 2 7 2 4 5 -1 5 3 14 56 4 0 14 5 13 -1 4 1 -2 0 0 1 2 -4 4 2 15 2 1 8 4 -2 1 11 11 0 7 12 4 7 2 0 7 3 1 19 11 3 0 1 1 1 5 2 3 3 0 -4 6 3 10 2 1 5 2 -1 15 -1 1 -4 11 12 2 4 3 5 36 3 4 1 4 4 5 3 7 7 10 1 1 -3 -3 0 1 6 4 4 3 1 5 -2 1 2 -1 4 14 4 4 4 -1 6 5 3 15 4 1 2 -2 0 1 1 21 -1 3 1 4 3 5 6 2 7 8 5 1 18 3 6 2 6 17 2 9 5 7 3 2 7 7 20 32 0 7 1 25 5 6 3 6 0 3 0 7 0 0 3 13 -2 2 2 3 0 6 0 2 7 3 2 2 5 1 6 5 13 -2 2 3 4 4 4 9 -1 2 5 2 9 2 2 -2 8 7 3 0 5 4 6 6 2 2 12 3 8 0 2 4 70 3 -1 2 1 6 2 5 0 5 1 10 2 2 2 2 7 14 3 1 3 1 11 3 4 1 0 5 2 11 1 3 11 5 10 3 3 3 2 0 2 1 3 -1 6 10 6 3 3 3 2 5 4 2 3 8 1 3 2 5 8 1 8 5 3 4 6 9 2 9 4 4 4 8 7 1 0 3 2 8 1 5 2 15 4 4 6 6 17 4 5 3 6 6 3 5 12 0 0 14 0 4 3 7 -2 5 3 14 0 10 2 2 7 4 13 2 2 2 3 2 1 0 22 9 15 0 4 3 4 2 2 8 0 8 6 10 5 3 1 3 0 1 9 6 5 2 4 5 0 3 3 0 3 6 13 7 -3 2 -2 2 2 4 7 -4 0 1 3 2 5 2 10 -3 4 2 4 4 9 7 5 9 3 8 12 -1 3 4 6 0 -1 1 8 6 -1 11 6 13 3 0 13 3 7 2 1 14 4 1 0 3 4 21 2 3 7 10 7 6 3 2 5 10 3 3 4 6 2 7 -1 3 4 9 2 2 1 3 -3 2 -1 -1 1 7 4 5 3 2 11 5 7 2 14 18 6 1 27 3 5 12 5 4 13 5 13 12 1 1 6 3 4 3 4 0 2 2 4 11 2 0 6 1 6 2 3 0 2 3 4 1 3 4 4 2 3 2 2'
$wsYDS.Range("C3").Value = '11 28 6 2 18 3 10 2 21 17 14 13 9 10 18 1 11 5 11 11 8 17 15 2 28 18 9 -1 14 16 12 2 12 15 4 8 5 9 20 14 11 12 9 1 14 7 14 17 20 14 26 8 14 19 6 5 11 19 19 18 9 20 16 9 16 8 8 8 10 5 7 11 5 4 -2 0 15 -1 4 7 -1 1 2 7 8 7 7 55 6 18 10 14 3 13 4 10 8 5 26 6 21 11 22 11 8 8 9 -7 10 15 19 38 9 4 3 29 10 8 11 20 8 12 2 0 8 7 5 11 6 4 5 10 11 5 18 5 3 13 6 8 22 10 8 9 18 15 13 12 8 16 14 6 10 1 15 16 7 3 6 11 40 59 3 11 30 18 9 19 3 10 3 7 14 3 9 2 8 9 15 12 8 1 9 8 4 16 3 20 25 8 18 21 7 11 11 12 6 0 7 6 14 32 45 8 6 33 7 15 3 6 8 8 68 8 5 8 7 10 11 4 8 8 5 22 19 15 6 1 6 5 16 8 16 5 9 13 7 6 12 9 16 11 1 5 1 10 2 4 13 7 14 4 5 18 16 1 5 6 24 6 6 6 4 12 18 6 21 -6 5 14 1 5 3 10 11 9 5 21 12 28 6 15 10 5 5 15 17 1 8 20 1 7 2 12 24 9 9 7 4 8 10 13 -3 7 41 11 11 3 5 8 17 4 6 17 10 22 9 2 7 18 15 22 5 2 12 18 20 8 22 6 17 8 19 -5 13 13 9 6 11 6 11 20 0 7 25 6 10 4 8 27 7 4 9 3 13 16 9 6 14 9 6 11 7 21 7 10 -3 5 20 7 10 8 3 16 3 10 5 3 5 7 15 13 7 6 1 9 5 25 -4 5 11 6 2 6 17 16 2 10 6 8 20 7 22 13 3 7 5 22 12 4 4 16 11 8 12 0 24 12 7 6 10 7 4 4 8 18 9 0 3 10 7 56 19 2 7 9 5 7 6 3 8 19 3 4 9 20 4 6 12 14 3 4 1 12 3 -1 10 26 4 6 14 4 1 1 28 25 -3 11 9 6 18 2 11 6 5 6 19 15 58 9 21 15 3 16 8 6 8 6 2 49 6 8 -2 11 21 7 15 24 11 30 15 4 8 14 18 22 35 7 30 16 6 28 10 15 9 4 4 -6 9 -4 16 25 4 4 5 18 5 6 13 4 2 9 1 3 5 5 12 21 5 12 6 13 7 5 1 7 7 19 10 3 9 12 7 0 15 5 6 6 7 11 9 17 9 3 11 12 9 18 8 4 9 10 10 15 14 24 2 6 6 20 6 2 8 3 7 8 10 20 7 2 15 9 5 9 20 24 22 11 2 3 12 3 19 9 17 8 3 10 8 15 3 13 7 10 8 6 6 6 7 4 35 7 12 16 3 13 9 17 14 14 18 4 5 8 3 21 5 3 14 0 19 2 0 13 4 21 6 10 8 7 14 11 9 5 11 1 4 12 8 6 -2 6 8 4 8 5 9 2 9 17 4 6 5 3 17 25 11 13 3 2 16 17'

$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B4").Value = '64 70 66 71 51 58 70 57 61 55 60 56 56 51 65 68 47 65 59 60 66 64 62 59 68 48 50 46 49 60 50 63 65 64 51 60 63 47 63 55 65 65 63 62 65 64 68 56 53 68 67 58 65 47 54 62 56 52 61 65 64 62 56 64 51 52 63 61 57 64'
$wsST.Range("B5").Value = '33 23 30 30 10 14 41 13 15 22 23 4 21 16 22 28 13 25 26 21 24 20 31 22 3 13 23 13 12 25 22 23 18 22 11 16 16 23 21 18 22 22 21 27 16 23 24 10 16 16 24 18 19 7 16 13 10 14 26 32 22 20 18 13 15 13 15 26 15 17'
$wsST.Range("D3").Value = '41 32 32 53 62 53 59 54 50 49 44 34 39 55 49 24 46 43 40 56 48 37 41 38 33 35 38 49 23 45 41 71 38 56 48 52 37 50 33 61 46 49 41 40 43 42 54 54 48 51 52 34 55 42 34 35 54 45 42 51 41 38 48 48 43 58 61 28 42 53 55 49 41 34 50 43 51 33 30 40 50 57 57 33 63 36 44 48 50 51 45 51 51 34 34 42 58 46 52 39 51 51 46 33 36 42 39 36 41 45 46 44 36'
$wsST.Range("D4").Value = '0 1 0 12 0 0 0 0 7 0 9 0 0 3 0 0 2 0 0 11 0 0 2 0 0 0 0 0 0 0 12 0 0 0 0 15 0 0 0 29 5 13 0 0 0 0 24 6 8 0 17 0 8 18 0 6 17 19 0 0 0 0 0 -2 0 13 11 0 0 0 15 9 0 0 0 9 0 0 0 17 0 8 28 0 11 0 10 0 0 4 0 0 0 0 0 0 0 5 0 0 9 0 4 0 2 3 9 0 0 7 0 0 0'
$wsST.Range("D5").Value = '8 20 11 0 11 0 0 15 0 9 0 6 0 0 0 14 5 0 0 0 0 11 0 6 0 1 0 10 16 20 0 8 0 0 0 0 0 0 0 0 0 0 0 0 0 12 0 0 6 0 0 20 0 19 0 0 2 11 8 0 0 0 0 0 0 19 17 0 26 0 7 0 0 0 0 0 0 2 12 6 0 4 0 0 0 12 14 0 0 0 9 0 0 0 15 11 0 -4 0'
